# Updates cryptos list price/volume figures (and the Aave/ordi row swap)
# to match the refreshed GitHub Actions data pull.
#
# Note: several "Price" values (e.g. 1.00, 316.04) look numeric, so a plain
# Range.Value assignment would auto-convert them and drop significant
# trailing zeros. A leading single-quote forces Excel to keep them as text,
# exactly like typing '1.00 into a cell. In a single-quoted PowerShell
# string that leading quote is written as '' (doubled), so e.g. '''1.00'
# assigns the literal text 1.00 (as text, not the number 1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '42.660.30'
$ws.Range("E2").Value = '  +0.01%  '

# Row 3
$ws.Range("D3").Value = '2.529.95'
$ws.Range("E3").Value = '  +0.55%  '

# Row 4
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.05%  '

# Row 5
$ws.Range("D5").Value = '''316.04'
$ws.Range("E5").Value = '  -0.44%  '

# Row 6
$ws.Range("D6").Value = '''96.03'
$ws.Range("E6").Value = '  +0.26%  '

# Row 7
$ws.Range("D7").Value = '''0.576'
$ws.Range("E7").Value = '  -0.49%  '

# Row 8
$ws.Range("E8").Value = '  +0.02%  '

# Row 9
$ws.Range("D9").Value = '''0.532'
$ws.Range("E9").Value = '  -0.85%  '

# Row 10
$ws.Range("D10").Value = '''35.69'
$ws.Range("E10").Value = '  -1.14%  '

# Row 11
$ws.Range("E11").Value = '  -0.43%  '

# Row 12
$ws.Range("D12").Value = '''7.53'
$ws.Range("E12").Value = '  -0.66%  '

# Row 13
$ws.Range("E13").Value = '  -2.45%  '

# Row 14
$ws.Range("D14").Value = '2.923.16'
$ws.Range("E14").Value = '  +0.74%  '

# Row 15
$ws.Range("D15").Value = '2.506.37'
$ws.Range("E15").Value = '  -0.22%  '

# Row 16
$ws.Range("D16").Value = '''15.10'
$ws.Range("E16").Value = '  -2.80%  '

# Row 17
$ws.Range("D17").Value = '''0.849'
$ws.Range("E17").Value = '  -0.93%  '

# Row 18
$ws.Range("D18").Value = '42.788.48'
$ws.Range("E18").Value = '  +0.28%  '

# Row 19
$ws.Range("E19").Value = '  +3.53%  '

# Row 20
$ws.Range("D20").Value = '''12.77'
$ws.Range("E20").Value = '  -2.43%  '

# Row 21
$ws.Range("D21").Value = '0.0₃0963'
$ws.Range("E21").Value = '  -0.82%  '

# Row 22
$ws.Range("D22").Value = '''69.55'
$ws.Range("E22").Value = '  -2.45%  '

# Row 23
$ws.Range("D23").Value = '''250.93'
$ws.Range("E23").Value = '  -0.30%  '

# Row 24
$ws.Range("D24").Value = '''2.95'
$ws.Range("E24").Value = '  -1.17%  '

# Row 25
$ws.Range("E25").Value = '  +1.05%  '

# Row 26
$ws.Range("D26").Value = '''26.40'
$ws.Range("E26").Value = '  -2.01%  '

# Row 27
$ws.Range("E27").Value = '  +0.05%  '

# Row 28
$ws.Range("D28").Value = '''2.40'
$ws.Range("E28").Value = '  +2.18%  '

# Row 29
$ws.Range("D29").Value = '''40.83'
$ws.Range("E29").Value = '  +4.89%  '

# Row 30
$ws.Range("D30").Value = '''10.39'
$ws.Range("E30").Value = '  +3.04%  '

# Row 31
$ws.Range("D31").Value = '''5.91'
$ws.Range("E31").Value = '  +0.02%  '

# Row 32
$ws.Range("D32").Value = '''157.39'
$ws.Range("E32").Value = '  +0.90%  '

# Row 33
$ws.Range("E33").Value = '  +2.80%  '

# Row 34
$ws.Range("E34").Value = '  +4.01%  '

# Row 35
$ws.Range("D35").Value = '''3.34'
$ws.Range("E35").Value = '  -0.17%  '

# Row 36
$ws.Range("D36").Value = '''18.86'
$ws.Range("E36").Value = '  -2.89%  '

# Row 37
$ws.Range("D37").Value = '''0.0782'
$ws.Range("E37").Value = '  -0.61%  '

# Row 38
$ws.Range("E38").Value = '  -1.19%  '

# Row 39
$ws.Range("E39").Value = '  -1.16%  '

# Row 40
$ws.Range("E40").Value = '  +8.24%  '

# Row 41
$ws.Range("D41").Value = '''22.49'
$ws.Range("E41").Value = '  -6.36%  '

# Row 42
$ws.Range("D42").Value = '''3.81'
$ws.Range("E42").Value = '  -1.09%  '

# Row 43
$ws.Range("D43").Value = '''0.0304'
$ws.Range("E43").Value = '  +1.18%  '

# Row 44
$ws.Range("E44").Value = '  +0.23%  '

# Row 45
$ws.Range("D45").Value = '2.033.54'
$ws.Range("E45").Value = '  -0.36%  '

# Row 46
$ws.Range("D46").Value = '''3.26'
$ws.Range("E46").Value = '  -3.49%  '

# Row 47
$ws.Range("D47").Value = '''9.07'
$ws.Range("E47").Value = '  +2.67%  '

# Row 48
$ws.Range("D48").Value = '''84.25'
$ws.Range("E48").Value = '  -0.29%  '

# Row 49
$ws.Range("B49").Value = 'ordi'
$ws.Range("C49").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D49").Value = '''75.45'
$ws.Range("E49").Value = '  +2.79%  '

# Row 50
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").Value = '''105.88'
$ws.Range("E50").Value = '  +3.87%  '

# Row 51
$ws.Range("D51").Value = '2.775.83'
$ws.Range("E51").Value = '  +0.68%  '
